$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.633.45"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "1.685.50"
$ws.Range("E3").Value = "  -1.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +1.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.38"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("E6").Value = "  +0.87%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3933"
$ws.Range("E7").Value = "  -0.61%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4036"
$ws.Range("E8").Value = "  -0.78%  "
$ws.Range("B9").Value = "BinanceUSD"
$ws.Range("C9").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.004"
$ws.Range("E9").Value = "  +1.11%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.480"
$ws.Range("E10").Value = "  -2.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.02"
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08790"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.201"
$ws.Range("E13").Value = "  -1.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.38"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.019"
$ws.Range("E15").Value = "  +7.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001309"
$ws.Range("E16").Value = "  -1.19%  "
$ws.Range("D17").Value = "1.691.78"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "99.49"
$ws.Range("E18").Value = "  -1.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06998"
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.39"
$ws.Range("E20").Value = "  -0.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.975"
$ws.Range("E21").Value = "  +3.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").Value = "  +1.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.25"
$ws.Range("E23").Value = "  +0.37%  "
$ws.Range("D24").Value = "24.616.18"
$ws.Range("E24").Value = "  -0.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.264"
$ws.Range("E25").Value = "  +9.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.366"
$ws.Range("E26").Value = "  +2.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.62"
$ws.Range("E27").Value = "  +0.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.52"
$ws.Range("E28").Value = "  +2.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.184"
$ws.Range("E29").Value = "  +1.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.95"
$ws.Range("E30").Value = "  +1.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.578"
$ws.Range("E31").Value = "  +1.81%  "
$ws.Range("D32").Value = "1.879.79"
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.08541"
$ws.Range("E33").Value = "  -1.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.056"
$ws.Range("E34").Value = "  -3.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.160"
$ws.Range("E35").Value = "  -3.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.10"
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2721"
$ws.Range("E37").Value = "  -0.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.893"
$ws.Range("E38").Value = "  -2.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.28"
$ws.Range("E39").Value = "  -3.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09160"
$ws.Range("E40").Value = "  +1.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02705"
$ws.Range("E41").Value = "  -2.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.457"
$ws.Range("E42").Value = "  -1.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7574"
$ws.Range("E43").Value = "  -1.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.98"
$ws.Range("E44").Value = "  +4.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.576"
$ws.Range("E45").Value = "  +4.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7108"
$ws.Range("E46").Value = "  -1.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.214"
$ws.Range("E47").Value = "  +1.33%  "
$ws.Range("E48").Value = "  +0.88%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "138.99"
$ws.Range("E49").Value = "  -1.75%  "
$ws.Range("B50").Value = "Flow"
$ws.Range("C50").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.310"
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("E51").Value = "  -0.76%  "
